# Edit script implementing the target diff:
#  1) Update the cached "datetimeFigureOut" field text (10/11/2018 -> 12/11/2018)
#     on the slide master and all 11 slide layouts' Date placeholders.
#  2) Flip + reposition/resize the "Straight Connector 18" connector on slide 1.
#  3) Reposition/resize the "TextBox 1" (the "X" mark) shape on slide 1.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholder text: slide master + every custom layout.
# ---------------------------------------------------------------------------
$m = $p.SlideMaster

for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $shp = $m.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq "10/11/2018") {
            $shp.TextFrame.TextRange.Text = "12/11/2018"
        }
    }
}

for ($l = 1; $l -le $m.CustomLayouts.Count; $l++) {
    $layout = $m.CustomLayouts.Item($l)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "10/11/2018") {
                $shp.TextFrame.TextRange.Text = "12/11/2018"
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) & 3) Shape geometry tweaks on slide 1.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

$connector = $s.Shapes.Item("Straight Connector 18")
$connector.HorizontalFlip = -1
$connector.Left = 427.8741149902344
$connector.Top = 182.4247283935547
$connector.Width = 0.8591338992118835
$connector.Height = 219.92913818359375

$crossMark = $s.Shapes.Item("TextBox 1")
$crossMark.Left = 409.4411315917969
$crossMark.Top = 379.9674072265625
$crossMark.Width = 32.9536247253418
$crossMark.Height = 46.0452766418457
